$d = $word.ActiveDocument

# The outline uses 9-level numbered/bulleted lists (see word/numbering.xml);
# LibreOffice-authored documents keep one "ListLabel N" character style per
# list level (1-9) alongside the existing "Bullets" character style. Add the
# nine missing ones now so the numbering infrastructure has the character
# styles it expects (same shape as the existing "Bullets" style: qFormat +
# a single complex-script font override of OpenSymbol).
1..9 | ForEach-Object {
    $styleName = "ListLabel " + $_
    $style = $d.Styles.Add($styleName, 2)   # wdStyleTypeCharacter
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}

Write-Output "Added ListLabel 1-9 character styles"
